# Weekly data refresh: a new record is prepended to the data table
# (which starts at row 2, the newest row always at row 103 in this
# sheet's history) and every subsequent row shifts down by one.
#
# Concretely: insert a new row at row 103 (pushing old rows 103..212
# down to 104..213) and populate the new row 103 with the latest
# week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 103:212 down to 104:213, leaving a blank
# row 103 for the new record.
$ws.Rows("103:103").Insert()

# Populate the newly inserted row with this week's data point.
$newRow = 103
$ws.Cells.Item($newRow, 1).Value = 10
$ws.Cells.Item($newRow, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($newRow, 3).Value = "La Araucanía"
$ws.Cells.Item($newRow, 4).Value = 44447
$ws.Cells.Item($newRow, 5).Value = 9
$ws.Cells.Item($newRow, 6).Value = 100112008
$ws.Cells.Item($newRow, 7).Value = "Coliflor"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 1250
$ws.Cells.Item($newRow, 11).Value = 800
$ws.Cells.Item($newRow, 12).Value = 800
$ws.Cells.Item($newRow, 13).Value = 800
$ws.Cells.Item($newRow, 14).Value = "$/unidad"
$ws.Cells.Item($newRow, 15).Value = "Región Metropolitana"
$ws.Cells.Item($newRow, 16).Value = 800
$ws.Cells.Item($newRow, 17).Value = 1
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
